$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the stray text value in E7 with the numeric time value used by the
# rest of that column (11:00 -> 0.45833333333333331), matching cells E1/E3/E5/E9.
$ws.Range("E7").Value = 0.45833333333333331

# Update the saved selection to E8.
$ws.Range("E8").Select()
